$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are treated as text so values such as
# "8.40" or "0.998" are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.600.67"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "3.279.24"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "577.15"
$ws.Range("E5").Value = "  +3.62%  "
$ws.Range("D6").Value = "183.37"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "3.275.02"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "0.568"
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("E10").Value = "  -5.33%  "
$ws.Range("E11").Value = "  -2.16%  "
$ws.Range("D12").Value = "46.08"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "3.817.82"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "8.40"
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "612.43"
$ws.Range("E16").Value = "  -3.25%  "
$ws.Range("D17").Value = "65.588.16"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "17.75"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "3.285.78"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "10.89"
$ws.Range("E21").Value = "  -4.06%  "
$ws.Range("D22").Value = "0.887"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "100.93"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").Value = "4.95"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "4.00"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "9.43"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("D29").Value = "30.76"
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("D30").Value = "8.42"
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("D31").Value = "6.42"
$ws.Range("D32").Value = "3.71"
$ws.Range("E32").Value = "  -9.12%  "
$ws.Range("D33").Value = "550.54"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("D35").Value = "3.783.56"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "55.90"
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("D39").Value = "0.127"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("D40").Value = "32.41"
$ws.Range("E40").Value = "  -4.21%  "
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("D42").Value = "3.37"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("E44").Value = "  -8.28%  "
$ws.Range("D45").Value = "0.330"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").Value = "3.03"
$ws.Range("E47").Value = "  -6.97%  "
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("D50").Value = "2.50"
$ws.Range("E50").Value = "  -3.73%  "
$ws.Range("D51").Value = "128.41"
$ws.Range("E51").Value = "  +4.47%  "
